$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.475.23'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.637.23'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("E5").Value = '  -0.26%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '303.83'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.13%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3784'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '51.69'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.30%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3613'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08158'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.39%  '

$ws.Range("E11").Value = '  -1.90%  '

$ws.Range("E12").Value = '  -0.13%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.45'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.34%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.448'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -3.21%  '

$ws.Range("E15").Value = '  +0.36%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001237'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.631.00'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.01%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06939'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.581'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.49'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -3.91%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.27%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.480.63'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.491'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.052'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.27%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.58%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '151.63'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.245'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.91'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.812.66'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.68%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.091'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +12.87%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.166'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.94%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.583'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.44'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.59%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02759'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.96%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2496'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.80%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08764'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.07081'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.973'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -5.24%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7037'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.61%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.347'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.32%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '12.17'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.61%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.54'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.51%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6530'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.74%  '

$ws.Range("E46").Value = '  -0.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.271'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.19%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.962'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.04%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07972'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.35%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.77'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.10%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.50%  '
